$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dados")

# Insert two new rows after row 2, both copies of row 2
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(4).Insert()

# Update row 3 values
$ws.Range("A3").Value = "Godzilla: King Of The Monsters 2"
$ws.Range("B3").Value = "Godzilla: King Of The Monsters 2"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "godzilla_king_otm_hd_net_evod_sub_ptbr2.ts"

# Update row 4 values
$ws.Range("A4").Value = "Godzilla: King Of The Monsters 3"
$ws.Range("B4").Value = "Godzilla: King Of The Monsters 3"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "godzilla_king_otm_hd_net_evod_sub_ptbr2.ts"

# Rebuild data validations so their ranges cover the new rows (3 and 4)
$ws.Cells.Validation.Delete()
$ws.Range("X2:X4").Validation.Add(3, 1, 1, "pais", "0")
$ws.Range("J2:J4").Validation.Add(3, 1, 1, "versao", "0")
$ws.Range("T2:U4").Validation.Add(3, 1, 1, "genero", "0")
$ws.Range("N2:N15").Validation.Add(3, 1, 1, "Categoria", "0")

# Update selection / view to match the target sheet view
$ws.Range("C11").Select()

Write-Host "done"
